$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 28 (shifts existing rows 28-47 down to 29-48,
# and extends the used range to A1:R48).
$ws.Rows.Item(28).Insert()

# Populate the newly inserted row 28 with the new record.
$ws.Range("A28").Value = 11
$ws.Range("B28").Value = "Vega Monumental Concepción"
$ws.Range("C28").Value = "Bíobío"
$ws.Range("D28").Value = 45161
$ws.Range("E28").Value = 8
$ws.Range("F28").Value = 100112022
$ws.Range("G28").Value = "Arveja Verde"
$ws.Range("H28").Value = "Perfection"
$ws.Range("I28").Value = "Primera"
$ws.Range("J28").Value = 50
$ws.Range("K28").Value = 32000
$ws.Range("L28").Value = 32000
$ws.Range("M28").Value = 32000
$ws.Range("N28").Value = '$/malla 25 kilos'
$ws.Range("O28").Value = "Provincia de Limarí"
$ws.Range("P28").Value = 1280
$ws.Range("Q28").Value = 25
$ws.Range("R28").Value = "Hortaliza"
